$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mock1")

# Update read-filtering counts (DADA2 option handling change)
$ws.Range("C2").Value = 8858
$ws.Range("D2").Value = 8859

$ws.Range("C4").Value = 1514
$ws.Range("D4").Value = 1514

$ws.Range("C6").Value = 930
$ws.Range("D6").Value = 930
